$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 8 for "Racialised" (shifts old row 8 "Ethnic.appearance" and below down by 1)
$ws.Rows.Item(8).Insert()
$ws.Range("A8").Value = "Racialised"

# Insert a new row at row 11 for "Ethnic.appearance.abridged"
# After the first insert: row8=Racialised, row9=Ethnic.appearance, row10=Ethnic.Appearance.original, row11=Found
$ws.Rows.Item(11).Insert()
$ws.Range("A11").Value = "Ethnic.appearance.abridged"
